$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2379
$ws.Range("I32").Value = 1623.6666
$ws.Range("J32").Value = 2832.2
$ws.Range("K32").Value = 1623.6666
$ws.Range("L32").Value = 2832.2
$ws.Range("M32").Value = -1297.6666
$ws.Range("N32").Value = -3484.2

$ws.Range("H33").Value = 4350.5
$ws.Range("I33").Value = 4350.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 4350.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -4121.5

$ws.Range("H38").Value = 149.5
$ws.Range("I38").Value = 149.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 448.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -76.5

$ws.Range("H51").Value = 5685.4287
$ws.Range("I51").Value = 4075
$ws.Range("J51").Value = 7832.6665
$ws.Range("K51").Value = 4075
$ws.Range("L51").Value = 7832.6665
$ws.Range("M51").Value = -3591
$ws.Range("N51").Value = -8800.666499999999

$ws.Range("H88").Value = 1883.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1883.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1883.5
$ws.Range("N88").Value = -2695.5

$ws.Range("H91").Value = 1883.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1883.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1883.5
$ws.Range("N91").Value = -4691.5

$ws.Range("H96").Value = 17392.834
$ws.Range("I96").Value = 34176
$ws.Range("J96").Value = 609.6667
$ws.Range("K96").Value = 102528
$ws.Range("L96").Value = 1829.0001
$ws.Range("M96").Value = -101155
$ws.Range("N96").Value = -4575.0001

$ws.Range("H131").Value = 18273.285
$ws.Range("I131").Value = 20783.6
$ws.Range("J131").Value = 11997.5
$ws.Range("K131").Value = 62350.8
$ws.Range("L131").Value = 35992.5
$ws.Range("M131").Value = -57310.8
$ws.Range("N131").Value = -46072.5

$ws.Range("H138").Value = 4613.8125
$ws.Range("I138").Value = 2528.25
$ws.Range("J138").Value = 5865.15
$ws.Range("K138").Value = 7584.75
$ws.Range("L138").Value = 17595.45
$ws.Range("M138").Value = -2444.75
$ws.Range("N138").Value = -27875.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2987.361
$ws.Range("I32").Value = 3085.0334
$ws.Range("J32").Value = 2499
$ws.Range("K32").Value = 3085.0334
$ws.Range("L32").Value = 2499
$ws.Range("M32").Value = -2798.0334
$ws.Range("N32").Value = -3073

$ws.Range("H97").Value = 690.4666999999999
$ws.Range("I97").Value = 528.4167
$ws.Range("J97").Value = 1338.6666
$ws.Range("K97").Value = 528.4167
$ws.Range("L97").Value = 1338.6666
$ws.Range("M97").Value = -32.41669999999999
$ws.Range("N97").Value = -2330.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2336
$ws.Range("I107").Value = 2082.3914
$ws.Range("J107").Value = 3065.125
$ws.Range("K107").Value = 2082.3914
$ws.Range("L107").Value = 3065.125
$ws.Range("M107").Value = -162.3914
$ws.Range("N107").Value = -6905.125

$ws.Range("H134").Value = 3748.6365
$ws.Range("I134").Value = 3748.6365
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11245.9095
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8710.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 170.6875
$ws.Range("I7").Value = 169.1
$ws.Range("J7").Value = 173.33333
$ws.Range("K7").Value = 169.1
$ws.Range("L7").Value = 173.33333
$ws.Range("M7").Value = -56.09999999999999
$ws.Range("N7").Value = -399.33333

$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -150
$ws.Range("N22").Value = -1700

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").ClearContents()
$ws.Range("N28").Value = 0

$ws.Range("H31").Value = 2016
$ws.Range("I31").Value = 2065.0557
$ws.Range("J31").Value = 1905.625
$ws.Range("K31").Value = 2065.0557
$ws.Range("L31").Value = 1905.625
$ws.Range("M31").Value = -1770.0557
$ws.Range("N31").Value = -2495.625

$ws.Range("H34").Value = 2016
$ws.Range("I34").Value = 2065.0557
$ws.Range("J34").Value = 1905.625
$ws.Range("K34").Value = 2065.0557
$ws.Range("L34").Value = 1905.625
$ws.Range("M34").Value = -1863.0557
$ws.Range("N34").Value = -2309.625

$ws.Range("H58").Value = 1294.4
$ws.Range("I58").Value = 1097.5
$ws.Range("J58").Value = 2082
$ws.Range("K58").Value = 1097.5
$ws.Range("L58").Value = 2082
$ws.Range("M58").Value = -894.5
$ws.Range("N58").Value = -2488

$ws.Range("H99").Value = 2265.5557
$ws.Range("I99").Value = 1865
$ws.Range("J99").Value = 3066.6667
$ws.Range("K99").Value = 1865
$ws.Range("L99").Value = 3066.6667
$ws.Range("M99").Value = -367
$ws.Range("N99").Value = -6062.6667

$ws.Range("H105").Value = 2860.9375
$ws.Range("I105").Value = 1559.375
$ws.Range("J105").Value = 4162.5
$ws.Range("K105").Value = 1559.375
$ws.Range("L105").Value = 4162.5
$ws.Range("M105").Value = 187.625
$ws.Range("N105").Value = -7656.5

$ws.Range("H107").Value = 1185
$ws.Range("I107").Value = 1094.6428
$ws.Range("J107").Value = 1395.8334
$ws.Range("K107").Value = 1094.6428
$ws.Range("L107").Value = 1395.8334
$ws.Range("M107").Value = 825.3571999999999
$ws.Range("N107").Value = -5235.8334

$ws.Range("H126").Value = 2265.5557
$ws.Range("I126").Value = 1865
$ws.Range("J126").Value = 3066.6667
$ws.Range("K126").Value = 5595
$ws.Range("L126").Value = 9200.000100000001
$ws.Range("M126").Value = -3125
$ws.Range("N126").Value = -14140.0001

$ws.Range("H136").Value = 1294.4
$ws.Range("I136").Value = 1097.5
$ws.Range("J136").Value = 2082
$ws.Range("K136").Value = 3292.5
$ws.Range("L136").Value = 6246
$ws.Range("M136").Value = -742.5
$ws.Range("N136").Value = -11346

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 374.4
$ws.Range("I7").Value = 344.5
$ws.Range("J7").Value = 394.33334
$ws.Range("K7").Value = 1033.5
$ws.Range("L7").Value = 1183.00002
$ws.Range("M7").Value = -921.5
$ws.Range("N7").Value = -1407.00002

$ws.Range("H12").Value = 811.9091
$ws.Range("I12").Value = 775.4286
$ws.Range("J12").Value = 875.75
$ws.Range("K12").Value = 2326.2858
$ws.Range("L12").Value = 2627.25
$ws.Range("M12").Value = -2153.2858
$ws.Range("N12").Value = -2973.25

$ws.Range("H74").Value = 13000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 39000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -41122

$ws.Range("H77").Value = 13000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 117000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -127608

$ws.Range("H129").Value = 574101.7
$ws.Range("I129").Value = 1455.6
$ws.Range("J129").Value = 1003586.3
$ws.Range("K129").Value = 4366.799999999999
$ws.Range("L129").Value = 3010758.9
$ws.Range("M129").Value = 633.2000000000007
$ws.Range("N129").Value = -3020758.9

$ws.Range("H138").Value = 8334969.5
$ws.Range("I138").Value = 11112791
$ws.Range("J138").Value = 1506
$ws.Range("K138").Value = 33338373
$ws.Range("L138").Value = 4518
$ws.Range("M138").Value = -33333233
$ws.Range("N138").Value = -14798

$ws.Range("H140").Value = 10924.0625
$ws.Range("I140").Value = 1148
$ws.Range("J140").Value = 14182.75
$ws.Range("K140").Value = 3444
$ws.Range("L140").Value = 42548.25
$ws.Range("M140").Value = 1736
$ws.Range("N140").Value = -52908.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3501966.8
$ws.Range("I7").Value = 5900
$ws.Range("J7").Value = 5250000
$ws.Range("K7").Value = 5900
$ws.Range("L7").Value = 5250000
$ws.Range("M7").Value = -5788
$ws.Range("N7").Value = -5250224

$ws.Range("H8").Value = 3501966.8
$ws.Range("I8").Value = 5900
$ws.Range("J8").Value = 5250000
$ws.Range("K8").Value = 5900
$ws.Range("L8").Value = 5250000
$ws.Range("M8").Value = -5761
$ws.Range("N8").Value = -5250278

$ws.Range("H11").Value = 10232538
$ws.Range("I11").Value = 13002300
$ws.Range("J11").Value = 1000000
$ws.Range("K11").Value = 13002300
$ws.Range("L11").Value = 1000000
$ws.Range("M11").Value = -13002161
$ws.Range("N11").Value = -1000278

$ws.Range("H107").Value = 1622.5862
$ws.Range("I107").Value = 1030.6471
$ws.Range("J107").Value = 2461.1667
$ws.Range("K107").Value = 1030.6471
$ws.Range("L107").Value = 2461.1667
$ws.Range("M107").Value = 889.3529000000001
$ws.Range("N107").Value = -6301.1667

$ws.Range("H132").Value = 2375.5454
$ws.Range("I132").Value = 1570.4445
$ws.Range("J132").Value = 5998.5
$ws.Range("K132").Value = 4711.333500000001
$ws.Range("L132").Value = 17995.5
$ws.Range("M132").Value = -2181.333500000001
$ws.Range("N132").Value = -23055.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1969.8
$ws.Range("I22").Value = 1949.6666

$ws.Range("H27").Value = 1969.8
$ws.Range("I27").Value = 1949.6666

$ws.Range("H46").Value = 3759.8
$ws.Range("I46").Value = 1519.6
$ws.Range("J46").Value = 6000
$ws.Range("K46").Value = 1519.6
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -1331.6
$ws.Range("N46").Value = -6376

$ws.Range("H68").Value = 2857.25
$ws.Range("I68").Value = 2850
$ws.Range("J68").Value = 2871.75
$ws.Range("K68").Value = 2850
$ws.Range("L68").Value = 2871.75
$ws.Range("M68").Value = -2101
$ws.Range("N68").Value = -4369.75

$ws.Range("H69").Value = 22387.666
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 22387.666
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 22387.666
$ws.Range("N69").Value = -24009.666

$ws.Range("H71").Value = 2857.25
$ws.Range("I71").Value = 2850
$ws.Range("J71").Value = 2871.75
$ws.Range("K71").Value = 14250
$ws.Range("L71").Value = 14358.75
$ws.Range("M71").Value = -10506
$ws.Range("N71").Value = -21846.75

$ws.Range("H72").Value = 22387.666
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 22387.666
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 67162.99800000001
$ws.Range("N72").Value = -75274.99800000001

$ws.Range("H132").Value = 4474.5
$ws.Range("I132").Value = 4533.8125
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 13601.4375
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -11071.4375
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 5000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -5228

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0

$ws.Range("H107").Value = 646.25
$ws.Range("I107").Value = 673.05884
$ws.Range("J107").Value = 604.8182
$ws.Range("K107").Value = 2019.17652
$ws.Range("L107").Value = 1814.4546
$ws.Range("M107").Value = -99.17651999999998
$ws.Range("N107").Value = -5654.4546

$ws.Range("H113").Value = 627.0714
$ws.Range("I113").Value = 531.4545000000001
$ws.Range("J113").Value = 977.6667
$ws.Range("K113").Value = 1594.3635
$ws.Range("L113").Value = 2933.0001
$ws.Range("M113").Value = 575.6364999999998
$ws.Range("N113").Value = -7273.0001

$ws.Range("H122").Value = 2414.5833
$ws.Range("I122").Value = 2179.5454
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6538.6362
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4088.6362
$ws.Range("N122").Value = -19900
